$d = $word.ActiveDocument

# The second paragraph currently reads:
#   "Prossime attività per settimana prossima:"
# built from the runs:
#   "Prossime attività" | " per settimana prossima" | <_GoBack bookmark> | ":"
#
# The edit turns "Prossime attività" into "Attività" (drop the leading
# "Prossime ", capitalise the initial letter) and leaves the _GoBack
# bookmark sitting right after the new initial "A", i.e. the final runs are:
#   "A" | <_GoBack bookmark> | "ttività" | " per settimana prossima" | ":"

$full = $d.Content.Text
$idx = $full.IndexOf("Prossime attività")
if ($idx -lt 0) {
    throw "Could not locate 'Prossime attività' in the document"
}

# Step 1 - relocate the _GoBack bookmark to sit between "Prossime a" and
# "ttività" (10 characters into the phrase). Re-adding a bookmark with the
# same name moves it there; doing this first splits the original run at
# that exact point while leaving every other run in the paragraph (the
# " per settimana prossima" run, the ":" run, ...) completely untouched.
$bookmarkPos = $idx + 10
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# Step 2 - "Prossime a" is now isolated in its own run immediately before
# the bookmark. Replace just that run's text with "A"; everything from the
# bookmark onward (the rest of "ttività", " per settimana prossima", ":")
# is left exactly as it was.
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("Prossime a")
if ($idx2 -lt 0) {
    throw "Could not locate 'Prossime a' in the document"
}
$r = $d.Range($idx2, $idx2 + 10)
$r.Text = "A"
